# Generate Report for handoff
# Fill in the actual "Latest Handoff Datetime" for the row whose
# "Latest Handoff File" is the 2a4cd2b0-... handoff file, on both the
# zh-cn and de-de locale report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-14 04:47:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-14 04:47:46"
